$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.780.91"
$ws.Range("E2").Value = "  -2.68%  "

$ws.Range("D3").Value = "1.746.38"
$ws.Range("E3").Value = "  -4.92%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "238.74"
$ws.Range("E5").Value = "  -8.87%  "

$ws.Range("E6").Value = "  -0.01%  "

$ws.Range("D7").Value = "0.5062"
$ws.Range("E7").Value = "  -6.00%  "

$ws.Range("D8").Value = "42.05"
$ws.Range("E8").Value = "  -6.28%  "

$ws.Range("D9").Value = "0.2731"
$ws.Range("E9").Value = "  -9.05%  "

$ws.Range("D10").Value = "0.06168"
$ws.Range("E10").Value = "  -10.98%  "

$ws.Range("D11").Value = "1.746.78"
$ws.Range("E11").Value = "  -4.98%  "

$ws.Range("D12").Value = "0.06934"
$ws.Range("E12").Value = "  -3.71%  "

$ws.Range("D13").Value = "15.55"
$ws.Range("E13").Value = "  -11.53%  "

$ws.Range("D14").Value = "4.524"
$ws.Range("E14").Value = "  -9.40%  "

$ws.Range("D15").Value = "0.6010"
$ws.Range("E15").Value = "  -18.19%  "

$ws.Range("D16").Value = "77.28"
$ws.Range("E16").Value = "  -13.36%  "

$ws.Range("D17").Value = "1.0000"
$ws.Range("E17").Value = "  -0.05%  "

$ws.Range("E18").Value = "  -0.01%  "

$ws.Range("D19").Value = "25.792.37"
$ws.Range("E19").Value = "  -2.75%  "

$ws.Range("D20").Value = "0.000006910"
$ws.Range("E20").Value = "  -12.55%  "

$ws.Range("E21").Value = "  -15.68%  "

$ws.Range("D22").Value = "1.970.78"
$ws.Range("E22").Value = "  -5.02%  "

$ws.Range("D23").Value = "4.068"
$ws.Range("E23").Value = "  -11.10%  "

$ws.Range("D24").Value = "5.261"
$ws.Range("E24").Value = "  -12.14%  "

$ws.Range("D25").Value = "8.190"
$ws.Range("E25").Value = "  -11.05%  "

$ws.Range("D26").Value = "137.73"
$ws.Range("E26").Value = "  -3.61%  "

$ws.Range("D27").Value = "1.471"
$ws.Range("E27").Value = "  -14.24%  "

$ws.Range("D28").Value = "1.815"
$ws.Range("E28").Value = "  -16.50%  "

$ws.Range("D29").Value = "15.00"
$ws.Range("E29").Value = "  -11.75%  "

$ws.Range("D30").Value = "103.86"
$ws.Range("E30").Value = "  -6.51%  "

$ws.Range("D31").Value = "0.08135"
$ws.Range("E31").Value = "  -8.15%  "

$ws.Range("D32").Value = "3.710"
$ws.Range("E32").Value = "  -12.46%  "

$ws.Range("D33").Value = "3.490"
$ws.Range("E33").Value = "  -13.77%  "

$ws.Range("D34").Value = "0.04528"
$ws.Range("E34").Value = "  -6.49%  "

$ws.Range("D35").Value = "0.9991"
$ws.Range("E35").Value = "  -0.03%  "

$ws.Range("D37").Value = "0.9835"
$ws.Range("E37").Value = "  -13.13%  "

$ws.Range("D38").Value = "0.6100"
$ws.Range("E38").Value = "  -16.28%  "

$ws.Range("D39").Value = "2.685"
$ws.Range("E39").Value = "  -13.20%  "

$ws.Range("D40").Value = "0.01553"
$ws.Range("E40").Value = "  -9.47%  "

$ws.Range("D41").Value = "1.937"
$ws.Range("E41").Value = "  -15.79%  "

$ws.Range("D42").Value = "1.0000"
$ws.Range("E42").Value = "  +0.01%  "

$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "102.03"
$ws.Range("E43").Value = "  -5.54%  "

$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.3847"
$ws.Range("E44").Value = "  -18.41%  "

$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "0.7397"
$ws.Range("E45").Value = "  -18.49%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "4.914"
$ws.Range("E46").Value = "  -16.83%  "

$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "0.05395"
$ws.Range("E47").Value = "  -6.45%  "

$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "0.1114"
$ws.Range("E48").Value = "  -11.00%  "

$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "5.974"
$ws.Range("E49").Value = "  -19.50%  "

$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "30.10"
$ws.Range("E50").Value = "  -13.58%  "

$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "52.50"
$ws.Range("E51").Value = "  -12.55%  "
